$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-28 20:58:12"
$ws.Range("E1:F1").ColumnWidth = 16.333333333333332

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-28 20:58:08"
$ws.Range("C1").ColumnWidth = 16.333333333333332

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-28 20:58:12"
$ws.Range("C1").ColumnWidth = 16.333333333333332
